{"js": "// Apply the \"Process\" -> \"process\"/\"method\" wording fixes described by the\n// commit. Each change is done as a small, targeted find & replace so that\n// unrelated occurrences of \"Process\" (e.g. \"Main Process\" headings,\n// \"Main_Process.py\" file name references) are left untouched.\n\nasync function replaceOnce(scopeRange, searchText, replacement, options) {\n  const opts = Object.assign({ matchCase: true, matchWholeWord: false }, options || {});\n  const results = scopeRange.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(searchText) +\n      \" but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) \"Main_Process.py is the Process,\" -> \"... is the process,\"\nawait replaceOnce(\n  body,\n  \"is the Process, which is the combination\",\n  \"is the process, which is the combination\"\n);\n\n// 2) \"Both of these Process deals\" -> \"Both of these process deals\"\nawait replaceOnce(\n  body,\n  \"Both of these Process deals\",\n  \"Both of these process deals\"\n);\n\n// 3) \"Primary Process related code\" -> \"Primary process-related code\"\nawait replaceOnce(\n  body,\n  \"Primary Process related code is written\",\n  \"Primary process-related code is written\"\n);\n\n// 4) \"this Process will perform the following steps\" -> \"... process ...\"\nawait replaceOnce(\n  body,\n  \"this Process will perform the following steps\",\n  \"this process will perform the following steps\"\n);\n\n// 5) \"flowchart of the Process for a known person.\" -> \"... process ...\"\nawait replaceOnce(\n  body,\n  \"flowchart of the Process for a known person.\",\n  \"flowchart of the process for a known person.\"\n);\n\n// 6) \"flowchart of the Process for unknown persons.\" -> \"... process ...\"\nawait replaceOnce(\n  body,\n  \"flowchart of the Process for unknown persons.\",\n  \"flowchart of the process for unknown persons.\"\n);\n\n// The sentence \"The primary Process will first call the Face Recognition\n// process ... this Process will call Speech Recognition.\" appears twice in\n// the document (once inline inside the intro paragraph, once as its own\n// paragraph later). The two copies need slightly different edits, so find\n// each owning paragraph individually before doing the scoped replace.\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet introParagraph = null; // contains \"...Processes. The primary Process...\"\nlet standaloneParagraph = null; // starts with \"The primary Process...\"\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"individual Processes. The primary Process will first call\") !== -1) {\n    introParagraph = paragraphs.items[i];\n  } else if (text.indexOf(\"The primary Process will first call\") === 0) {\n    standaloneParagraph = paragraphs.items[i];\n  }\n}\n\nif (!introParagraph) {\n  throw new Error(\"Could not locate the intro paragraph with 'The primary Process'.\");\n}\nif (!standaloneParagraph) {\n  throw new Error(\"Could not locate the standalone 'The primary Process' paragraph.\");\n}\n\n// 7) Intro paragraph: \"The primary Process\" -> \"The primary method\"\nawait replaceOnce(introParagraph, \"The primary Process will\", \"The primary method will\");\n\n// 8) Intro paragraph: \"this Process will call Speech Recognition.\" -> \"this process ...\"\nawait replaceOnce(\n  introParagraph,\n  \"this Process will call Speech Recognition.\",\n  \"this process will call Speech Recognition.\"\n);\n\n// 9) Standalone paragraph: \"The primary Process\" -> \"The primary process\"\nawait replaceOnce(standaloneParagraph, \"The primary Process will\", \"The primary process will\");\n\n// 10) Standalone paragraph: \"this Process will call Speech Recognition.\" -> \"this process ...\"\nawait replaceOnce(\n  standaloneParagraph,\n  \"this Process will call Speech Recognition.\",\n  \"this process will call Speech Recognition.\"\n);\n", "ps1": "# Apply the \"Process\" -> \"process\"/\"method\" wording fixes described by the\n# commit. Each change is done as a small, targeted find & replace, scoped to\n# a specific Range so that unrelated occurrences of \"Process\" (e.g. \"Main\n# Process\" headings, \"Main_Process.py\" file name references) are left\n# untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-InRange($range, [string]$findText, [string]$replaceText) {\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0          # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $result = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $result) {\n        throw \"Find/replace failed for '$findText'\"\n    }\n}\n\n# 1) \"Main_Process.py is the Process,\" -> \"... is the process,\"\nReplace-InRange $d.Content \"is the Process, which is the combination\" \"is the process, which is the combination\"\n\n# 2) \"Both of these Process deals\" -> \"Both of these process deals\"\nReplace-InRange $d.Content \"Both of these Process deals\" \"Both of these process deals\"\n\n# 3) \"Primary Process related code\" -> \"Primary process-related code\"\nReplace-InRange $d.Content \"Primary Process related code is written\" \"Primary process-related code is written\"\n\n# 4) \"this Process will perform the following steps\" -> \"... process ...\"\nReplace-InRange $d.Content \"this Process will perform the following steps\" \"this process will perform the following steps\"\n\n# 5) \"flowchart of the Process for a known person.\" -> \"... process ...\"\nReplace-InRange $d.Content \"flowchart of the Process for a known person.\" \"flowchart of the process for a known person.\"\n\n# 6) \"flowchart of the Process for unknown persons.\" -> \"... process ...\"\nReplace-InRange $d.Content \"flowchart of the Process for unknown persons.\" \"flowchart of the process for unknown persons.\"\n\n# The sentence \"The primary Process will first call the Face Recognition\n# process ... this Process will call Speech Recognition.\" appears twice in\n# the document (once inline inside the intro paragraph, once as its own\n# paragraph later). The two copies need slightly different edits, so find\n# each owning paragraph individually before doing the scoped replace.\n\n$introParagraph = $null\n$standaloneParagraph = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*individual Processes. The primary Process will first call*\") {\n        $introParagraph = $p\n    } elseif ($t -like \"The primary Process will first call*\") {\n        $standaloneParagraph = $p\n    }\n}\n\nif ($null -eq $introParagraph) {\n    throw \"Could not locate the intro paragraph with 'The primary Process'.\"\n}\nif ($null -eq $standaloneParagraph) {\n    throw \"Could not locate the standalone 'The primary Process' paragraph.\"\n}\n\n# 7) Intro paragraph: \"The primary Process\" -> \"The primary method\"\nReplace-InRange $introParagraph.Range() \"primary Process will\" \"primary method will\"\n\n# 8) Intro paragraph: \"this Process will call Speech Recognition.\" -> \"this process ...\"\nReplace-InRange $introParagraph.Range() \"this Process will call Speech Recognition.\" \"this process will call Speech Recognition.\"\n\n# 9) Standalone paragraph: \"The primary Process\" -> \"The primary process\"\nReplace-InRange $standaloneParagraph.Range() \"primary Process will\" \"primary process will\"\n\n# 10) Standalone paragraph: \"this Process will call Speech Recognition.\" -> \"this process ...\"\nReplace-InRange $standaloneParagraph.Range() \"this Process will call Speech Recognition.\" \"this process will call Speech Recognition.\"\n"}
